$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tabla")

# Add two new header columns after Z1 (Wb y), matching the existing
# header style (bold/centered/bordered, same as B1:Z1).
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("AA1:AB1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("AA1").Value = "Exp Constant"
$ws.Range("AB1").Value = "Exp Constant [dB]"
